# Regression script update for 1081, 898 and 518 stories
#  - bump the LiveSLR build number shown in the Application/Version table
#  - leave the active selection on the Version cell for that row (B2)
#  - restore the workbook window to its normal (non-maximized) size/position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2 holds "Copyright @ 2022 Cytel Inc. LiveSLR 4.0.0.0 - Build #49237"
# -> bump the build number to #49753
$ws.Range("B2").Value = "Copyright @ 2022 Cytel Inc. LiveSLR 4.0.0.0 - Build #49753"

# Move the active cell/selection from C2 to B2
$ws.Range("B2").Select()

# Restore the workbook window to a normal, non-maximized size/position
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 5760
$excel.ActiveWindow.Top = 3432
$excel.ActiveWindow.Width = 17280
$excel.ActiveWindow.Height = 9072
